# Commit: "Thres update to multiple area abs sum"
# On the SP_REP_DELAY sheet, a new "otfs-sbl" column (C) is inserted that
# keeps a copy of the previous threshold values (old column B), while
# column B is updated with new threshold values for a handful of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SP_REP_DELAY")

# Header for the new column C ("otfs-sbl")
$ws.Range("C1").Value = "otfs-sbl"

# Column C gets a snapshot of the ORIGINAL column B values (rows 2-17)
$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 15
$ws.Range("C4").Value = 16
$ws.Range("C5").Value = 15
$ws.Range("C6").Value = 14
$ws.Range("C7").Value = 37

# C8 mirrors B8, including its scientific number format (style index 3)
$ws.Range("C8").NumberFormat = $ws.Range("B8").NumberFormat
$ws.Range("C8").Value = 10000

$ws.Range("C9").Value = 4
$ws.Range("C10").Value = 2
$ws.Range("C11").Value = $ws.Range("B11").Value2
$ws.Range("C12").Value = $ws.Range("B12").Value2
$ws.Range("C13").Value = 6
$ws.Range("C14").Value = 3
$ws.Range("C15").Value = 5
$ws.Range("C16").Value = 240
$ws.Range("C17").Value = 240

# Column B gets new updated threshold values for these rows
$ws.Range("B4").Value = 9
$ws.Range("B5").Value = 7
$ws.Range("B14").Value = 2
$ws.Range("B15").Value = 3
$ws.Range("B16").Value = 63
$ws.Range("B17").Value = 63

# Selection now targets B17 only (was A16:B17)
$ws.Range("B17").Select()
